$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update results for Steel (Iron & steel / Hydrogen)
$ws.Range("B3").Value = 6162.199529381775

# Minor recalculated precision updates for Non-metallic minerals column
$ws.Range("D6").Value = 3862.877934945198
$ws.Range("D8").Value = 428.1943212310456
